# Update the "Sources" annexe table with the newer source-refresh dates
# (28/10/2020 -> 21/01/2021, IMF WEO 2018 -> 2019 vintage, etc.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sources")

# FMI (IMF) WEO database row
$ws.Range("D4").Value = "Mis à jour le 21/01/2021"

# IMF WEO dataset URL moved from the 2018 vintage to the 2019 vintage
$ws.Range("E5").Value = "https://www.imf.org/external/pubs/ft/weo/2019/01/weodata/index.aspx"

# OECD gender-index database refresh date
$ws.Range("D6").Value = "Décembre 2019"

# UNDP human development index vintage year
$ws.Range("D13").Value = "2019"

# World Bank WDI indicators refresh date
$ws.Range("D17").Value = "Mis à jour le 16/12/2020"

# Wittgenstein Centre human capital data vintage year
$ws.Range("D22").Value = "2019"

# UN World Urbanization Prospects vintage year (title + year column)
$ws.Range("B25").Value = "World Urbanization Prospects: the 2019 revision, ONU"
$ws.Range("D25").Value = "2019"

# The E5 hyperlink keeps pointing at the original (2018) address while the
# cell text now shows the 2019 URL, so Excel records the old address as the
# cached display text for the link.
$ws.Hyperlinks.Add($ws.Range("E5"), "https://www.imf.org/external/pubs/ft/weo/2018/01/weodata/index.aspx", "", "", "https://www.imf.org/external/pubs/ft/weo/2018/01/weodata/index.aspx")
$ws.Range("E5").Value = "https://www.imf.org/external/pubs/ft/weo/2019/01/weodata/index.aspx"

# Restored window size from the author's last save
$excel.ActiveWindow.Height = 587
